$d = $word.ActiveDocument

$pairs = @(
    @("54×38=2052", "92×87=8004"),
    @("72×44=3168", "16×30=480"),
    @("64×99=6336", "15×59=885"),
    @("33×55=1815", "38×63=2394"),
    @("60×67=4020", "47×29=1363"),
    @("55×20=1100", "60×58=3480"),
    @("85×57=4845", "96×77=7392"),
    @("63×75=4725", "65×74=4810"),
    @("32×87=2784", "79×28=2212"),
    @("13×90=1170", "39×58=2262"),
    @("18×64=1152", "61×38=2318"),
    @("57×75=4275", "32×68=2176"),
    @("12×20=240",  "32×23=736"),
    @("42×53=2226", "18×79=1422"),
    @("38×47=1786", "13×55=715"),
    @("65×92=5980", "24×41=984"),
    @("93×61=5673", "86×66=5676"),
    @("12×68=816",  "66×78=5148"),
    @("89×16=1424", "65×21=1365"),
    @("22×91=2002", "36×69=2484"),
    @("55×79=4345", "60×21=1260"),
    @("16×56=896",  "80×40=3200"),
    @("61×83=5063", "15×65=975"),
    @("95×92=8740", "17×47=799"),
    @("76×38=2888", "27×34=918")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
